$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the input values that drive all downstream formulas.
$ws.Range("C2").Value = 20842
$ws.Range("C3").Value = 21500
$ws.Range("C4").Value = 5863
$ws.Range("C5").Value = 1775
$ws.Range("C6").Value = 909600
$ws.Range("C7").Value = 9279000
$ws.Range("C8").Value = 9099000
$ws.Range("C9").Value = 5100000
$ws.Range("C10").Value = 3200000

# Remove bold formatting across the whole sheet.
$ws.Cells.Font.Bold = $false
